$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete row 9 (the old "NB" / index-7 row is removed entirely) ---
$ws.Rows.Item(9).Delete()

# --- 2. Extend header formatting from G1 into H1:L1 (copy format only) ---
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Header row text (B1:L1) ---
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# --- 4. Algorithm name column (B2:B8); "CART" renamed to "DTREE", "NB" row removed ---
$ws.Range("B2").Value = "LR"
$ws.Range("B3").Value = "LDA"
$ws.Range("B4").Value = "KNN"
$ws.Range("B5").Value = "DTREE"
$ws.Range("B6").Value = "RTREE"
$ws.Range("B7").Value = "XTREE"
$ws.Range("B8").Value = "SVM"

# --- 5. Numeric data, columns C:L for rows 2-8 (mean/std pairs per horizon) ---
$data = @{
    2 = @(0.9090633192976881, 0.007655612320712047, 0.8944535291078293, 0.004195076109444465, 0.882236415073543, 0.01079074962255617, 0.8722299793778999, 0.01990929111612589, 0.8608850808732044, 0.01811944010881063)
    3 = @(0.9141068827626387, 0.008389764514451451, 0.9015338886333828, 0.008642836735566218, 0.8888877322327409, 0.01071340011758976, 0.8777206135074757, 0.01816703426875145, 0.8627898427779662, 0.01175707424204308)
    4 = @(0.8976756441831519, 0.00860070140328904, 0.8959715195802886, 0.0147702220321703, 0.8866104123547205, 0.01274066556905293, 0.8909710288136565, 0.01402472819133107, 0.8839503449835991, 0.01400964459002119)
    5 = @(0.8874290935091762, 0.01234268693575006, 0.8850125765808734, 0.0095418996474835, 0.8850394351707835, 0.01489291970865616, 0.8775333677035, 0.01458242444907673, 0.8856187082909173, 0.01298724687587197)
    6 = @(0.9093882577262253, 0.006807428383615113, 0.8976553051595211, 0.007481318467676062, 0.8857396543911431, 0.008423643936893786, 0.8739330927421664, 0.01699277216502831, 0.8635058251329035, 0.01674456319474749)
    7 = @(0.9173618283414104, 0.009721805478364327, 0.9087843017016709, 0.009730136102976808, 0.8990358590621288, 0.01108614190671006, 0.8994919659735349, 0.01578586165891406, 0.8970320099536252, 0.009927898663925837)
    8 = @(0.908899393554196, 0.005620607179130398, 0.9057511597140602, 0.00910680787134438, 0.9007874787209289, 0.008782040317103978, 0.8979775305035229, 0.01592426866106061, 0.8891884402216945, 0.01338668063846961)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 3 + $i   # column C = 3
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
